$wb = $excel.ActiveWorkbook

# 1. Remove the empty Plan2 and Plan3 tabs, keep only Plan1
$wb.Worksheets("Plan2").Delete() | Out-Null
$wb.Worksheets("Plan3").Delete() | Out-Null

$ws = $wb.Worksheets("Plan1")

# 2. Fill in the previously-empty simulation results for table 1 (rows 2-5)
$ws.Range("B2").Value = 173
$ws.Range("C2").Value = 123
$ws.Range("B3").Value = 1810
$ws.Range("C3").Value = 1109
$ws.Range("B4").Value = 18201
$ws.Range("C4").Value = 11103
$ws.Range("D4").Value = 130
$ws.Range("B5").Value = 175804
$ws.Range("C5").Value = 122945
$ws.Range("D5").Value = 646

# 3. Note introducing the refactored (second) simulation table
$ws.Range("A11").Value = "DEPOIS DA REFATORAÇÃO PARA CAIR MENOS AVIÕES"

# 4. Second table header (row 13) - same column titles as row 1
$ws.Range("A13").Value = "numeroMaxIteracao"
$ws.Range("B13").Value = "avioesDecolados"
$ws.Range("C13").Value = "avioesPousados"
$ws.Range("D13").Value = "avioesCaidos"

# 5. Second table data (rows 14-17)
$ws.Range("A14").Value = 100
$ws.Range("B14").Value = 198
$ws.Range("C14").Value = 102
$ws.Range("D14").Value = 0

$ws.Range("A15").Value = 1000
$ws.Range("B15").Value = 1977
$ws.Range("C15").Value = 937
$ws.Range("D15").Value = 0

$ws.Range("A16").Value = 10000
$ws.Range("B16").Value = 19846
$ws.Range("C16").Value = 9791
$ws.Range("D16").Value = 0

$ws.Range("A17").Value = 100000
$ws.Range("B17").Value = 198199
$ws.Range("C17").Value = 101259
$ws.Range("D17").Value = 0

# 6. Style the second header row (row 13) the same way as the first (row 1):
#    bold text, with the data columns (B-D) centered
$ws.Range("A13:D13").Font.Bold = $true
$ws.Range("B13:D13").HorizontalAlignment = -4108
$ws.Range("B13:D13").VerticalAlignment = -4108

# 7. Put a thin box border around both data tables
$ws.Range("A1:D5").Borders.LineStyle = 1
$ws.Range("A13:D17").Borders.LineStyle = 1

# 8. Thousands-separator number format for the larger counters
$ws.Range("A3:C5").NumberFormat = "#,##0"
$ws.Range("A15:B17,C17").NumberFormat = "#,##0"

# 9. Leave the selection on the newly added second table, like the source edit
$ws.Range("A13:D17").Select()
